# Applies the "Add registration info / Add symposium and young stats session"
# edit to sessions.xlsx:
#  1. The "Featured Sessions" row that used to hold the STRATOS satellite
#     symposium (Title/Organizer) is replaced with the new
#     "Young Statisticians Sessions and Panel Discussion" session.
#  2. A new "Satellite Symposium" worksheet is appended (after "Featured
#     Sessions") holding the STRATOS row that moved out of Featured Sessions,
#     plus a blank trailing row, and becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

$featured = $wb.Worksheets.Item("Featured Sessions")

# Row 6 on "Featured Sessions" currently holds the STRATOS satellite
# symposium entry -- replace it with the new Young Statisticians session.
$featured.Cells.Item(6, 1).Value = "Young Statisticians Sessions and Panel Discussion"
$featured.Cells.Item(6, 2).Value = "Andrea Berghold, Stefanie Peschel"

# Add the new "Satellite Symposium" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$satellite = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$satellite.Name = "Satellite Symposium"

$satellite.Cells.Item(1, 1).Value = "Title"
$satellite.Cells.Item(1, 2).Value = "Organizer"
$satellite.Cells.Item(2, 1).Value = "Ten years of the STRengthening Analytical Thinking for Observational Studies (STRATOS) initiative – progress and looking to the future"
$satellite.Cells.Item(2, 2).Value = "Ruth Keogh, Willi Sauerbrei"

$satellite.Activate()
